$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = "2025-02-27"

$data = @(
    ,@("service", 1, 1)
    ,@("gouvernement", 1, 2)
    ,@("autorisation", 1, 1)
    ,@("service", 2, 17)
    ,@("gouvernement", 2, 1)
    ,@("autorisation", 2, 1)
    ,@("service", 4, 12)
    ,@("gouvernement", 4, 1)
    ,@("service", 5, 12)
    ,@("service", 6, 7)
    ,@("service", 7, 1)
    ,@("gouvernement", 7, 1)
    ,@("autorisation", 7, 3)
    ,@("autorisation", 8, 1)
    ,@("gouvernement", 11, 1)
    ,@("gouvernement", 12, 3)
    ,@("gouvernement", 13, 1)
    ,@("gouvernement", 14, 2)
    ,@("autorisation", 14, 3)
    ,@("service", 14, 1)
    ,@("service", 15, 1)
    ,@("gouvernement", 15, 3)
    ,@("service", 16, 6)
    ,@("service", 17, 7)
    ,@("service", 18, 9)
    ,@("autorisation", 18, 3)
    ,@("service", 30, 1)
    ,@("service", 36, 1)
    ,@("autorisation", 36, 1)
    ,@("service", 39, 1)
    ,@("service", 42, 2)
    ,@("service", 53, 1)
    ,@("gouvernement", 57, 2)
    ,@("gouvernement", 59, 4)
    ,@("service", 64, 8)
    ,@("service", 65, 13)
    ,@("service", 66, 5)
    ,@("service", 67, 3)
    ,@("service", 68, 2)
    ,@("service", 69, 1)
    ,@("autorisation", 69, 1)
    ,@("service", 72, 2)
    ,@("autorisation", 74, 4)
    ,@("autorisation", 75, 1)
    ,@("service", 78, 8)
    ,@("service", 79, 4)
)

$startRow = 42
$endRow = $startRow + $data.Count - 1

# Force column A to be stored as text (not auto-converted to a date serial number),
# matching the inlineStr/shared-string representation used for the date values,
# then strip the resulting number-format style so no style index is left on the cell.
$dateRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $term = $data[$i][0]
    $page = $data[$i][1]
    $occ = $data[$i][2]
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $term
    $ws.Cells.Item($row, 3).Value = $page
    $ws.Cells.Item($row, 4).Value = $occ
}

$dateRange.ClearFormats()
